# logBook.xlsx update — "updated till 3rd july 230pm"
# Adds two new log entries (rows 40 & 41) plus a trailing blank formatted
# row (42), matches the existing row 38/39 formatting, and refreshes the
# view state / total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Seed rows 40 and 41 by copying the formatting of the two rows
#        immediately above them (37/38 -> 40/41 pattern already used by the
#        sheet), then overwrite with the new values/formulas below. This
#        guarantees the new cells reuse the same style indices (date/time/
#        wrap-text formats) as the rest of the table.
$ws.Range("A38:G38").Copy($ws.Range("A40:G40"))
$ws.Range("A39:G39").Copy($ws.Range("A41:G41"))

# --- 2. Row 40 : Sno 39, 3-Jul-2022, 10:00 - 11:30, Code
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 44745
$ws.Range("C40").Value = 0.41666666666666669
$ws.Range("D40").Value = 0.47916666666666669
$ws.Range("E40").Formula = "=D40-C40"
$ws.Range("F40").Value = "Code"
$ws.Range("G40").Value = "1. Block diagrams for PPM modules`n2. Block diagrams for PSPNet architecture and Aux loss"

# --- 3. Row 41 : Sno 40, 3-Jul-2022, 13:15 - 14:30, Code
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 44745
$ws.Range("C41").Value = 0.55208333333333337
$ws.Range("D41").Value = 0.60416666666666663
$ws.Range("E41").Formula = "=D41-C41"
$ws.Range("F41").Value = "Code"
$ws.Range("G41").Value = "1. Formatted and uploaded PSPNet_starter nb`n2. PSPNet_resnet50_aux nb completed"

# Row 40/41 match the 2-line wrapped description row height (30pt) used
# elsewhere in the sheet for shorter two-line entries.
$ws.Rows.Item(40).RowHeight = 30
$ws.Rows.Item(41).RowHeight = 30

# --- 4. Row 42 : trailing blank row that keeps the same column formatting
#        (date/time/duration/wrap styles) but carries no values, same as
#        the blank rows that already exist between data and the totals row.
$ws.Range("B41:E41").Copy()
$ws.Range("B42:E42").PasteSpecial(-4122)
$ws.Range("G41").Copy()
$ws.Range("G42").PasteSpecial(-4122)
$ws.Range("B42:E42").ClearContents()
$ws.Range("G42").ClearContents()

# --- 5. Refresh the Total Hours formula result (same formula, now spans
#        the two new rows of data).
$ws.Range("E50").Formula = "=SUM(E2:E49)"

# --- 6. Update the saved view state to match where the sheet was left:
#        scrolled down to row 38 and the active cell on D49.
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 1
$ws.Range("D49").Select()

$wb.Application.CalculateFull()
